$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reformat the "bday" column (I) and the "telephone_number" column (R) ---
# The bday column is switching from a Date number format to plain Text, and
# picks up the bordered / centered header look that telephone_number (R) used
# to have. telephone_number loses that special formatting and goes back to
# the default look.

# Clear the old bordered/centered + number-format styling from column R
# (this also drops the column-level style hint Excel had stamped on col R).
$ws.Columns("R").ClearFormats()

# Give the bday header (I1) the bordered, centered look, now stored as Text.
$ws.Range("I1").Borders.LineStyle = 1
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").NumberFormat = "@"

# The bday data cells (I2:I4) switch from Date format to Text format, so the
# new non-numeric-month values ("2024-26-01", etc.) are stored verbatim.
$ws.Range("I2:I4").NumberFormat = "@"

# --- Header text tweak ---
$ws.Range("I1").Value = "bday (YEAR-DATE-MONTH)"

# --- Data edits (seeder fix) ---
$ws.Range("G2").Value = "ruby"
$ws.Range("I2").Value = "2024-07-02"

$ws.Range("G3").Value = "humility"
$ws.Range("I3").Value = "2024-05-03"

$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "lapu-lapu"
$ws.Range("I4").Value = "2024-26-01"

# --- Sidenav / view fix: scroll the sheet over and move the active selection ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("I7").Select()
